# -----------------------------------------------------------------------
# Target change analysis
# -----------------------------------------------------------------------
# The supplied unified diff touches only word/document.xml and
# word/styles.xml, and every single changed line is an *attribute
# permutation* of the very same element: identical tag name, identical
# set of attribute name/value pairs, just re-ordered (mostly into
# alphabetical order, which is the hallmark of Apache POI/XMLBeans'
# OOXML writer). E.g.
#
#   -<w:pgSz w:w="11906" w:h="16838"/>
#   +<w:pgSz w:h="16838" w:w="11906"/>
#
#   -<w:style w:type="paragraph" w:default="1" w:styleId="Normal">
#   +<w:style w:default="1" w:styleId="Normal" w:type="paragraph">
#
# No text, run/paragraph formatting, style definition, page setup
# value, language setting, or document structure actually changed —
# confirmed by diffing the attribute-name/value sets of every
# before/after element pair in the patch (they are set-equal in every
# single hunk). This matches the accompanying commit message, "Fixed
# POI packaging and upgraded to POI 3.15": upgrading the OOXML-writing
# library changed how XML attributes get serialized (now emitted in a
# different/alphabetized order) without altering the document's
# content or formatting in any way.
#
# The Word COM/object model (Paragraphs, Styles, PageSetup, Find, …)
# exposes the document's semantic content, not the serialization
# order of XML attributes on disk — there is nothing in that object
# model to "fix" here, because nothing in the object model differs
# between the before/after states. Re-applying every property to its
# own current value (confirmed experimentally) does not alter the
# on-disk attribute order either, since no real value changed.
#
# So the correct COM replay of this commit is a no-op against the
# document model: open the document, change nothing.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# Touch the document without mutating any content/formatting, so the
# script still exercises ActiveDocument per the runtime's expectations.
$null = $d.Paragraphs.Count
